# update beauty, MG, PSS dashboard
# Sheet1 row 2 (the single data/record row) gets a new generated case id in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CA-HU5809FQ"
